$wb = $excel.ActiveWorkbook

# --- Sheet "Means" ---
$wsMeans = $wb.Worksheets.Item("Means")

$wsMeans.Range("B9").Value = 23
$wsMeans.Range("C9").Value = 29
$wsMeans.Range("F9").Value = 47
$wsMeans.Range("G9").Value = 51

$wsMeans.Range("B10").Value = 0.27
$wsMeans.Range("C10").Value = 0.36
$wsMeans.Range("D10").Value = 0.5
$wsMeans.Range("E10").Value = 0.5
$wsMeans.Range("F10").Value = 0.5
$wsMeans.Range("G10").Value = 0.51

# --- Sheet "Standard Deviations" ---
$wsSD = $wb.Worksheets.Item("Standard Deviations")

$wsSD.Range("B9").Value = 7.2
$wsSD.Range("C9").Value = 5.4
$wsSD.Range("F9").Value = 5.2
$wsSD.Range("G9").Value = 12

$wsSD.Range("B10").Value = 0.094
$wsSD.Range("C10").Value = 0.058
$wsSD.Range("F10").Value = 0
$wsSD.Range("G10").Value = 0.077
